# "Add date formats, composite set"
#
# style.xlsx before this change has 6 fills / 4 cellXfs. This edit:
#   1. Recolors the fill used by C3's style (fillId 5) from FF22FF00 to 99BB66.
#   2. Adds a new composite style (solid fill + date numFmt) used by A4
#      (white fill, numFmtId 16 "d-mmm").
#   3. Adds a second new composite style used by A5
#      (orange fill, numFmtId 14 "mm-dd-yy").
#   A4/A5 both get the serial date value 39904.166666666664 (2009-04-01 04:00).
#
# Excel's Interior.Color is a BGR-packed long (0xBBGGRR), i.e. the reverse
# byte order of a normal "RRGGBB" hex string.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function RgbColor([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

$dateSerial = 39904.166666666664

# 1. Recolor fillId 5 (currently FF22FF00) -> 99BB66. The only cell using
#    that style is C3, so recoloring its Interior achieves the fill edit.
$ws.Range("C3").Interior.Color = RgbColor 0x99 0xBB 0x66

# 2. A4: new fill (white, FFFFFF) + built-in date format "d-mmm" (numFmtId 16)
$ws.Range("A4").Value = $dateSerial
$ws.Range("A4").Interior.Color = RgbColor 0xFF 0xFF 0xFF
$ws.Range("A4").NumberFormat = "d-mmm"

# 3. A5: new fill (orange, FFAA00) + built-in date format "mm-dd-yy" (numFmtId 14)
$ws.Range("A5").Value = $dateSerial
$ws.Range("A5").Interior.Color = RgbColor 0xFF 0xAA 0x00
$ws.Range("A5").NumberFormat = "mm-dd-yy"
